# SF-Pending-Bookings-List-Template.xlsx
# Rearrange the booking-list columns (Booking ID / TimeSlot / Age moved to the
# front, Product / Services regrouped, etc.) and refresh the view selection.
#
# The column order/content move is done the same way a person would do it in
# Excel: cut the whole column and insert it before its new home, letting the
# other columns shift across. Doing it as a short sequence of cut/insert
# operations (instead of rewriting every cell) automatically carries each
# column's formatting (style) along with its content, matching how Excel
# itself performs a column move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rearrange columns -------------------------------------------------
# Each step: cut the column currently holding a given header, insert it
# (shifting everything at/after the destination to the right) at its final
# position. Because earlier moves change where things live, the source
# letters below are computed for the state of the sheet *at that step*.

$ws.Columns("D").Cut()     # Booking ID      -> position A
$ws.Columns("A").Insert()

$ws.Columns("K").Cut()     # Product         -> position B
$ws.Columns("B").Insert()

$ws.Columns("E").Cut()     # Booking Date    -> position C (after A/B insertions)
$ws.Columns("C").Insert()

$ws.Columns("Q").Cut()     # Booking TimeSlot -> position D
$ws.Columns("D").Insert()

$ws.Columns("Q").Cut()     # Booking Age     -> position E
$ws.Columns("E").Insert()

$ws.Columns("N").Cut()     # Engineer Name   -> position M
$ws.Columns("M").Insert()

# --- 2) Re-apply the column widths shown in the final layout -------------
# Excel's ColumnWidth is expressed in "characters" of the Normal style font;
# the stored OOXML width is ColumnWidth + 5/6 (this sheet's Normal font has
# a 6px max-digit-width), so subtract 5/6 from each target width before
# assigning it.
$ws.Columns("A").ColumnWidth = 20.9440104166667
$ws.Columns("B").ColumnWidth = 18.1666666666667
$ws.Columns("C").ColumnWidth = 19.9440104166667
$ws.Columns("D").ColumnWidth = 24.7213541666667
$ws.Columns("E").ColumnWidth = 20.2760416666667
$ws.Columns("F").ColumnWidth = 21.8307291666667
$ws.Columns("G").ColumnWidth = 27.3854166666667
$ws.Columns("H").ColumnWidth = 15.0533854166667
$ws.Columns("I").ColumnWidth = 16.6080729166667
$ws.Columns("J").ColumnWidth = 16.6080729166667
$ws.Columns("K").ColumnWidth = 19.7213541666667
$ws.Columns("L").ColumnWidth = 33.7213541666667
$ws.Columns("M").ColumnWidth = 20.0533854166667
$ws.Columns("N").ColumnWidth = 18.4986979166667
$ws.Columns("O").ColumnWidth = 23.6080729166667
$ws.Columns("P").ColumnWidth = 17.8307291666667
$ws.Columns("Q").ColumnWidth = 19.0533854166667
$ws.Columns("R").ColumnWidth = 15.1666666666667

# --- 3) Update the visible selection --------------------------------------
$ws.Range("U13").Select()
